{"js": "// 1) Remove the existing \"_GoBack\" bookmark (it currently sits right after\n//    the page-break run, near the top of the document).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Every literal tab character run (`<w:tab/>`) in the lab's CREATE TABLE\n//    listings gets replaced by four literal space characters\n//    (`<w:t xml:space=\"preserve\">    </w:t>`). `Body.search` treats \"\\t\" as\n//    plain text and matches the tab runs one-for-one.\nconst body = context.document.body;\nconst tabs = body.search(\"\\t\", { matchCase: true });\ntabs.load(\"items\");\nawait context.sync();\n\nconst items = tabs.items;\nfor (let i = 0; i < items.length; i++) {\n  items[i].insertText(\"    \", \"Replace\");\n}\nawait context.sync();\n\n// 3) The \"_GoBack\" bookmark (last-edit marker) now belongs right after the\n//    very last tab-turned-spaces run (the last edit the author made).\nconst lastRange = items[items.length - 1];\nconst collapsed = lastRange.getRange(\"End\");\nawait context.sync();\ncollapsed.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) The existing \"_GoBack\" bookmark (last-edit marker) currently sits right\n#    after the page-break run near the top of the document; drop it, it will\n#    be re-created at the author's actual last edit location below.\n$d.Bookmarks(\"_GoBack\").Delete()\n\n# 2) Every literal tab character (`<w:tab/>`) in the lab's CREATE TABLE\n#    listings is retyped as four literal spaces\n#    (`<w:t xml:space=\"preserve\">    </w:t>`). Walk all tab hits with\n#    Find/Execute and retype each one in place.\n$rng = $d.Content\n$count = 0\nwhile ($rng.Find.Execute([ref]\"^t\")) {\n    $count = $count + 1\n    $rng.Text = \"    \"\n    $rng.Collapse(0)  # wdCollapseEnd - resume searching right after the fix\n\n    # 3) \"_GoBack\" belongs right after the very last tab-turned-spaces run\n    #    (the last place the author edited).\n    if ($count -eq 6) {\n        $d.Bookmarks.Add(\"_GoBack\", $rng)\n    }\n}\n"}
